$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New row of data (row 9) appended below the existing GILD bag-of-words rows.
$row = 9

# Copy the date cell's format (A8) down to A9 so the new timestamp keeps the
# same date/time number format as the rest of column A.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial($xlPasteFormats)

$ws.Cells.Item($row, 1).Value = 42611.883634259262
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 44
$ws.Cells.Item($row, 5).Value = 66
$ws.Cells.Item($row, 6).Value = 33
$ws.Cells.Item($row, 7).Value = 14019
$ws.Cells.Item($row, 8).Value = 28257
$ws.Cells.Item($row, 9).Value = 3210
$ws.Cells.Item($row, 10).Value = 389
$ws.Cells.Item($row, 11).Value = 339
$ws.Cells.Item($row, 12).Value = 12
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Bag"
